$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Perguntas com cores dadas iguais:" summary block (columns L/M), added
# alongside the existing question table.

# Header label (L3) - no special formatting.
$ws.Range("L3").Value = "Perguntas com cores dadas iguais:"

# Row 4: questions #3 + #12 share color #80FF00
$ws.Range("L4").Value = "#80FF00"
$ws.Range("L4").Interior.Color = 65408
$ws.Range("L4").HorizontalAlignment = -4108
$ws.Range("L4").VerticalAlignment = -4108
$ws.Range("L4").Borders.Item(10).Weight = -4138
$ws.Range("M4").Value = "#3 + #12"

# Row 5: questions #6 + #11 share color #FF8000
$ws.Range("L5").Value = "#FF8000"
$ws.Range("L5").Interior.Color = 33023
$ws.Range("L5").HorizontalAlignment = -4108
$ws.Range("L5").VerticalAlignment = -4108
$ws.Range("L5").Borders.Item(10).Weight = -4138
$ws.Range("M5").Value = "#6 + #11"

# Row 6: questions #9 + #15 share color #00FF80
$ws.Range("L6").Value = "#00FF80"
$ws.Range("L6").Interior.Color = 8453888
$ws.Range("L6").HorizontalAlignment = -4108
$ws.Range("L6").VerticalAlignment = -4108
$ws.Range("L6").Borders.Item(10).Weight = -4138
$ws.Range("M6").Value = "#9 + #15"

# Row 7: questions #10 + #13 share color #0080FF
$ws.Range("L7").Value = "#0080FF"
$ws.Range("L7").Interior.Color = 16744448
$ws.Range("L7").HorizontalAlignment = -4108
$ws.Range("L7").VerticalAlignment = -4108
$ws.Range("L7").Borders.Item(10).Weight = -4138
$ws.Range("M7").Value = "#10 + #13"

# Scroll back to the top and select the cell next to the new block, matching
# the saved view state (no more frozen topLeftCell at A10).
$ws.Range("A1").Select()
$ws.Range("M9").Select()
